# Rimuovi PEC Comando dalla risposta, perche devono passare tramite SUAP/SUE/SUA
#
# 1) Collapse "... questa comunicazione, tramite [indirizzo PEC] certificata:"
#    down to "... questa comunicazione." (drop the "tramite PEC" clause).
# 2) Delete the whole paragraph that holds the $PEC_COMANDO MERGEFIELD
#    (the <...> placeholder paragraph right after it).
# 3) Normal style: overflowPunct goes from true to false.

$d = $word.ActiveDocument

# --- 1) Trim the "Il richiedente e invitato a regolarizzare ..." sentence ---
$find = $d.Content.Find
$find.Execute(", tramite l’indirizzo di posta elettronica certificata:", $false, $false, $false, $false, $false, $true, 1, $false, ".", 2) | Out-Null

# --- 2) Remove the paragraph containing the $PEC_COMANDO mergefield result ---
$paras = $d.Paragraphs
$target = $null
for ($i = 1; $i -le $paras.Count; $i++) {
    $p = $paras.Item($i)
    if ($p.Range.Text -match "PEC_COMANDO" -and $i -gt 10) {
        $target = $p
    }
}
if ($target -ne $null) {
    $target.Range.Delete()
}

# --- 3) Normal style paragraph format: overflowPunct true -> false ---
$style = $d.Styles("Normal")
$style.ParagraphFormat.HangingPunctuation = $false
